$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 96.64286
$ws.Range("I5").Value = 62.75
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 62.75
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = 52.25
$ws.Range("N5").Value = -530

$ws.Range("H32").Value = 1010
$ws.Range("I32").Value = 1010
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1010
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -684
$ws.Range("N32").Value = ""

$ws.Range("H58").Value = 2292.182
$ws.Range("I58").Value = 604.6667
$ws.Range("J58").Value = 2925
$ws.Range("K58").Value = 1814.0001
$ws.Range("L58").Value = 8775
$ws.Range("M58").Value = -1664.0001
$ws.Range("N58").Value = -9075

$ws.Range("H92").Value = 1566.6666
$ws.Range("I92").Value = 4000
$ws.Range("J92").Value = 350
$ws.Range("K92").Value = 4000
$ws.Range("L92").Value = 350
$ws.Range("M92").Value = -2752

$ws.Range("H100").Value = 1303
$ws.Range("I100").Value = 1303
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1303
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -762

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = ""

$ws.Range("H125").Value = 1693.3334
$ws.Range("I125").Value = 1620.1428
$ws.Range("J125").Value = 1949.5
$ws.Range("K125").Value = 14581.2852
$ws.Range("L125").Value = 17545.5
$ws.Range("M125").Value = -12121.2852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 21823.334
$ws.Range("I31").Value = 21823.334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 21823.334
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -21529.334

$ws.Range("H32").Value = 2328756
$ws.Range("I32").Value = 2800.2
$ws.Range("J32").Value = 7696346.5
$ws.Range("K32").Value = 2800.2
$ws.Range("L32").Value = 7696346.5
$ws.Range("M32").Value = -2513.2
$ws.Range("N32").Value = -7696920.5

$ws.Range("H45").Value = 2807.8125
$ws.Range("I45").Value = 2184.4
$ws.Range("J45").Value = 3846.8333
$ws.Range("K45").Value = 2184.4
$ws.Range("L45").Value = 3846.8333
$ws.Range("M45").Value = -1807.4
$ws.Range("N45").Value = -4600.8333

$ws.Range("H63").Value = 7965.5
$ws.Range("I63").Value = 1898
$ws.Range("J63").Value = 10999.25
$ws.Range("K63").Value = 1898
$ws.Range("L63").Value = 10999.25
$ws.Range("M63").Value = -1212

$ws.Range("H66").Value = 7965.5
$ws.Range("I66").Value = 1898
$ws.Range("J66").Value = 10999.25
$ws.Range("K66").Value = 9490
$ws.Range("L66").Value = 54996.25
$ws.Range("M66").Value = -6058

$ws.Range("H132").Value = 3840.7778
$ws.Range("I132").Value = 3537.2942
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 10611.8826
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -8081.882599999999
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1147
$ws.Range("I64").Value = 1147
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1147
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -922
$ws.Range("N64").Value = ""

$ws.Range("H67").Value = 1147
$ws.Range("I67").Value = 1147
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 1147
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -367
$ws.Range("N67").Value = ""

$ws.Range("H80").Value = 828.53845
$ws.Range("I80").Value = 604.25
$ws.Range("J80").Value = 1187.4
$ws.Range("K80").Value = 604.25
$ws.Range("L80").Value = 1187.4
$ws.Range("M80").Value = 393.75
$ws.Range("N80").Value = -3183.4

$ws.Range("H83").Value = 828.53845
$ws.Range("I83").Value = 604.25
$ws.Range("J83").Value = 1187.4
$ws.Range("K83").Value = 3021.25
$ws.Range("L83").Value = 5937
$ws.Range("M83").Value = 1970.75
$ws.Range("N83").Value = -15921

$ws.Range("H86").Value = 5931.684
$ws.Range("I86").Value = 4251
$ws.Range("J86").Value = 7799.1113
$ws.Range("K86").Value = 4251
$ws.Range("L86").Value = 7799.1113
$ws.Range("M86").Value = -3128

$ws.Range("H89").Value = 5931.684
$ws.Range("I89").Value = 4251
$ws.Range("J89").Value = 7799.1113
$ws.Range("K89").Value = 21255
$ws.Range("L89").Value = 38995.5565
$ws.Range("M89").Value = -15639

$ws.Range("H94").Value = 1471.4546
$ws.Range("I94").Value = 1479.9474
$ws.Range("J94").Value = 1417.6666
$ws.Range("K94").Value = 1479.9474
$ws.Range("L94").Value = 1417.6666
$ws.Range("M94").Value = -1028.9474
$ws.Range("N94").Value = -2319.6666

$ws.Range("H99").Value = 1566.6666
$ws.Range("I99").Value = 1472.25
$ws.Range("J99").Value = 1755.5
$ws.Range("K99").Value = 1472.25
$ws.Range("L99").Value = 1755.5
$ws.Range("M99").Value = 25.75

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = ""

$ws.Range("H134").Value = 1381.9333
$ws.Range("I134").Value = 1381.9333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4145.7999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1610.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -187

$ws.Range("H7").Value = 3232.4856
$ws.Range("I7").Value = 6380.3125
$ws.Range("J7").Value = 581.6842
$ws.Range("K7").Value = 6380.3125
$ws.Range("L7").Value = 581.6842
$ws.Range("M7").Value = -6267.3125

$ws.Range("H22").Value = 34.666668
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 34.666668
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 34.666668
$ws.Range("N22").Value = -734.666668

$ws.Range("H62").Value = 872.5
$ws.Range("I62").Value = 745
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 745
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = -121
$ws.Range("N62").Value = -2248

$ws.Range("H65").Value = 872.5
$ws.Range("I65").Value = 745
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 3725
$ws.Range("L65").Value = 5000
$ws.Range("M65").Value = -605
$ws.Range("N65").Value = -11240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 8176.143
$ws.Range("I55").Value = 4945
$ws.Range("J55").Value = 9468.6
$ws.Range("K55").Value = 14835
$ws.Range("L55").Value = 28405.8
$ws.Range("M55").Value = -14658
$ws.Range("N55").Value = -28759.8

$ws.Range("H68").Value = 1873.2307
$ws.Range("I68").Value = 1716.8334
$ws.Range("J68").Value = 2007.2858
$ws.Range("K68").Value = 5150.5002
$ws.Range("L68").Value = 6021.857400000001
$ws.Range("M68").Value = -4339.5002
$ws.Range("N68").Value = -7643.857400000001

$ws.Range("H71").Value = 1873.2307
$ws.Range("I71").Value = 1716.8334
$ws.Range("J71").Value = 2007.2858
$ws.Range("K71").Value = 15451.5006
$ws.Range("L71").Value = 18065.5722
$ws.Range("M71").Value = -11395.5006
$ws.Range("N71").Value = -26177.5722

$ws.Range("H107").Value = 326.66666
$ws.Range("I107").Value = 330.66666
$ws.Range("J107").Value = 318.66666
$ws.Range("K107").Value = 991.9999799999999
$ws.Range("L107").Value = 955.9999799999999
$ws.Range("M107").Value = 928.0000200000001
$ws.Range("N107").Value = -4795.99998

$ws.Range("H131").Value = 1665.9412
$ws.Range("I131").Value = 975.25
$ws.Range("J131").Value = 2279.889
$ws.Range("K131").Value = 2925.75
$ws.Range("L131").Value = 6839.667
$ws.Range("M131").Value = 2114.25

$ws.Range("H132").Value = 1835.1177
$ws.Range("I132").Value = 1789.7
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 16107.3
$ws.Range("L132").Value = 17100
$ws.Range("M132").Value = -13577.3
$ws.Range("N132").Value = -22160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 9500
$ws.Range("I99").Value = 9500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 9500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -7254

$ws.Range("H140").Value = 39042
$ws.Range("I140").Value = 39697
$ws.Range("J140").Value = 38714.5
$ws.Range("K140").Value = 39697
$ws.Range("L140").Value = 38714.5
$ws.Range("M140").Value = -34517
$ws.Range("N140").Value = -49074.5

$ws.Range("H141").Value = 64999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 64999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 64999
$ws.Range("N141").Value = -75359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 31994
$ws.Range("I2").Value = 31994
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 31994
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -31882

$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 2250
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 2250
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1955
$ws.Range("N22").Value = -3590

$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 2250
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 2250
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -2143
$ws.Range("N27").Value = -3214

$ws.Range("H40").Value = 2794.8572
$ws.Range("I40").Value = 2794.8572
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2794.8572
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2658.8572

$ws.Range("H46").Value = 4050.375
$ws.Range("I46").Value = 2900.5
$ws.Range("J46").Value = 7500
$ws.Range("K46").Value = 2900.5
$ws.Range("L46").Value = 7500
$ws.Range("M46").Value = -2712.5
$ws.Range("N46").Value = -7876

$ws.Range("H136").Value = 3251.75
$ws.Range("I136").Value = 2127.625
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 6382.875
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -3832.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""

$ws.Range("H62").Value = 11999.941
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 12124.9375
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 12124.9375
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -13372.9375

$ws.Range("H65").Value = 11999.941
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 12124.9375
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 60624.6875
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -66864.6875

$ws.Range("H81").Value = 10400
$ws.Range("I81").Value = 9250
$ws.Range("J81").Value = 15000
$ws.Range("K81").Value = 18500
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = -17439
$ws.Range("N81").Value = -32122

$ws.Range("H84").Value = 10400
$ws.Range("I84").Value = 9250
$ws.Range("J84").Value = 15000
$ws.Range("K84").Value = 92500
$ws.Range("L84").Value = 150000
$ws.Range("M84").Value = -87196
$ws.Range("N84").Value = -160608

$ws.Range("H96").Value = 1387.5
$ws.Range("I96").Value = 1500
$ws.Range("J96").Value = 1275
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 1275
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -4021

$ws.Range("H126").Value = 4899.9287
$ws.Range("I126").Value = 3250.0625
$ws.Range("J126").Value = 7099.75
$ws.Range("K126").Value = 9750.1875
$ws.Range("L126").Value = 21299.25
$ws.Range("M126").Value = -7280.1875
$ws.Range("N126").Value = -26239.25

$ws.Range("H136").Value = 3914.7693
$ws.Range("I136").Value = 2570.5715
$ws.Range("J136").Value = 5483
$ws.Range("K136").Value = 7711.7145
$ws.Range("L136").Value = 16449
$ws.Range("M136").Value = -5161.7145
$ws.Range("N136").Value = -21549
